{"js": "// Update the date line and replace every arithmetic expression in the\n// practice-sheet table with the new values from the target revision.\n// (20 rows x 5 columns of \"a+b=\" / \"a-b=\" style cells.)\n\nconst newDate = \"2022-11-26 Saturday\";\n\nconst newGrid = [\n  [\"54-34=\", \"26+65=\", \"72+21=\", \"45+50=\", \"43+25=\"],\n  [\"99-18=\", \"66+31=\", \"7+1=\", \"89-84=\", \"66-26=\"],\n  [\"10-7=\", \"30+59=\", \"55-8=\", \"38-25=\", \"34+30=\"],\n  [\"59-43=\", \"18+70=\", \"79-48=\", \"93-27=\", \"69-32=\"],\n  [\"10-0=\", \"47+4=\", \"14+24=\", \"58+1=\", \"61-52=\"],\n  [\"50+14=\", \"7+22=\", \"43-19=\", \"66+13=\", \"54+15=\"],\n  [\"24-9=\", \"9+60=\", \"49+42=\", \"30+26=\", \"93-52=\"],\n  [\"20+50=\", \"8+4=\", \"45-9=\", \"36+57=\", \"76-26=\"],\n  [\"21+4=\", \"79-30=\", \"61+1=\", \"13-1=\", \"72-49=\"],\n  [\"97-0=\", \"69-18=\", \"40-22=\", \"73-31=\", \"27+66=\"],\n  [\"31-22=\", \"68-26=\", \"30+4=\", \"54+6=\", \"19+16=\"],\n  [\"32+6=\", \"3+57=\", \"55+26=\", \"4+45=\", \"96-82=\"],\n  [\"11+74=\", \"52-29=\", \"7+52=\", \"47+13=\", \"45-2=\"],\n  [\"82-81=\", \"88-83=\", \"33+11=\", \"35+29=\", \"55+38=\"],\n  [\"17+61=\", \"77+20=\", \"96-95=\", \"84+2=\", \"53+38=\"],\n  [\"28-8=\", \"71-65=\", \"25+18=\", \"73+3=\", \"9+67=\"],\n  [\"5+3=\", \"30+49=\", \"70-10=\", \"87+7=\", \"36-26=\"],\n  [\"47-2=\", \"80-44=\", \"43+7=\", \"14+10=\", \"15+56=\"],\n  [\"27+44=\", \"83-1=\", \"64-0=\", \"82-58=\", \"65-27=\"],\n  [\"81-51=\", \"79-60=\", \"76-2=\", \"71+22=\", \"57-38=\"],\n];\n\nconst body = context.document.body;\n\n// 1) Update the title paragraph with the date.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.getRange().insertText(newDate, Word.InsertLocation.replace);\n\n// 2) Update every cell in the practice table, preserving cell formatting.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newGrid;\n\nawait context.sync();\n", "ps1": "# Update the date line and replace every arithmetic expression in the\n# practice-sheet table with the new values from the target revision.\n# (20 rows x 5 columns of \"a+b=\" / \"a-b=\" style cells.)\n\n$d = $word.ActiveDocument\n\n$newDate = \"2022-11-26 Saturday\"\n\n$newGrid = @(\n    @(\"54-34=\", \"26+65=\", \"72+21=\", \"45+50=\", \"43+25=\"),\n    @(\"99-18=\", \"66+31=\", \"7+1=\", \"89-84=\", \"66-26=\"),\n    @(\"10-7=\", \"30+59=\", \"55-8=\", \"38-25=\", \"34+30=\"),\n    @(\"59-43=\", \"18+70=\", \"79-48=\", \"93-27=\", \"69-32=\"),\n    @(\"10-0=\", \"47+4=\", \"14+24=\", \"58+1=\", \"61-52=\"),\n    @(\"50+14=\", \"7+22=\", \"43-19=\", \"66+13=\", \"54+15=\"),\n    @(\"24-9=\", \"9+60=\", \"49+42=\", \"30+26=\", \"93-52=\"),\n    @(\"20+50=\", \"8+4=\", \"45-9=\", \"36+57=\", \"76-26=\"),\n    @(\"21+4=\", \"79-30=\", \"61+1=\", \"13-1=\", \"72-49=\"),\n    @(\"97-0=\", \"69-18=\", \"40-22=\", \"73-31=\", \"27+66=\"),\n    @(\"31-22=\", \"68-26=\", \"30+4=\", \"54+6=\", \"19+16=\"),\n    @(\"32+6=\", \"3+57=\", \"55+26=\", \"4+45=\", \"96-82=\"),\n    @(\"11+74=\", \"52-29=\", \"7+52=\", \"47+13=\", \"45-2=\"),\n    @(\"82-81=\", \"88-83=\", \"33+11=\", \"35+29=\", \"55+38=\"),\n    @(\"17+61=\", \"77+20=\", \"96-95=\", \"84+2=\", \"53+38=\"),\n    @(\"28-8=\", \"71-65=\", \"25+18=\", \"73+3=\", \"9+67=\"),\n    @(\"5+3=\", \"30+49=\", \"70-10=\", \"87+7=\", \"36-26=\"),\n    @(\"47-2=\", \"80-44=\", \"43+7=\", \"14+10=\", \"15+56=\"),\n    @(\"27+44=\", \"83-1=\", \"64-0=\", \"82-58=\", \"65-27=\"),\n    @(\"81-51=\", \"79-60=\", \"76-2=\", \"71+22=\", \"57-38=\")\n)\n\n# 1) Update the title paragraph with the date.\n$d.Paragraphs.Item(1).Range.Text = $newDate\n\n# 2) Update every cell in the practice table, preserving cell formatting.\n$table = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newGrid.Count; $r++) {\n    $row = $newGrid[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
